$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell labels (shared strings) per the commit diff ---
$ws.Range("C6").Value = "Obesity grade, acute COVID-19"
$ws.Range("C9").Value = "Days hospitalized, acute COVID-19"
$ws.Range("C16").Value = "Metabolic disorders"
$ws.Range("C22").Value = "Immune deficiency"
$ws.Range("C23").Value = "Weight loss, acute COVID-19"
$ws.Range("C24").Value = "Dyspnea, acute COVID-19"
$ws.Range("C25").Value = "Cough, acute COVID-19"
$ws.Range("C26").Value = "Fever, acute COVID-19"
$ws.Range("C27").Value = "Night sweat, acute COVID-19"
$ws.Range("C28").Value = "Pain, acute COVID-19"
$ws.Range("C29").Value = "GI symptoms, acute COVID-19"
$ws.Range("C30").Value = "Anosmia, acute COVID-19"
$ws.Range("C31").Value = "ECOG, acute COVID-19"
$ws.Range("C32").Value = "Impaired performance, acute COVID-19"
$ws.Range("C33").Value = "Sleep disorders, acute COVID-19"
$ws.Range("C34").Value = "Anti-infectives, acute COVID-19"
$ws.Range("C35").Value = "Anti-platelet, acute COVID-19"
$ws.Range("C36").Value = "Anti-coagulatives, acute COVID-19"
$ws.Range("C37").Value = "Immunosuppression, acute COVID-19"
$ws.Range("C39").Value = "# Self-reported symptoms, acute COVID-19"
$ws.Range("C40").Value = "Anti-S1/S2 IgG, 60-day visit"
$ws.Range("C41").Value = "Hb, 60-day visit"
$ws.Range("C42").Value = "Anemia, 60-day visit"
$ws.Range("C43").Value = "Ferritin, 60-day visit"
$ws.Range("C44").Value = "Elevated ferritin, 60-day visit"
$ws.Range("C45").Value = "sTFR, 60-day visit"
$ws.Range("C46").Value = "Hepcidin, 60-day visit"
$ws.Range("C47").Value = "Elevated NTproBNP, 60-day visit"
$ws.Range("C48").Value = "Elevated D-dimer, 60-day visit"
$ws.Range("C49").Value = "Elevated CRP, 60-day visit"
$ws.Range("C50").Value = "Elevated IL-6, 60-day visit"
$ws.Range("C51").Value = "Ferritin Index, 60-day visit"
$ws.Range("C52").Value = "Iron deficiency, 60-day visit"
$ws.Range("C53").Value = "Age over 65"
$ws.Range("C54").Value = "Hospitalized >7 days, acute COVID-19"
$ws.Range("C56").Value = ">3 comorbidities"
$ws.Range("C57").Value = "Overweight or obesity"
$ws.Range("C58").Value = ">6 symptoms, acute COVID-19"
$ws.Range("C59").Value = "Persistent symptoms, 60-day visit"
$ws.Range("E59").Value = "> 0 symptoms at 180-day visit"
$ws.Range("C60").Value = "Anti-S1/S2 IgG Q1, 60-day visit"
$ws.Range("C61").Value = "Anti-S1/S2 IgG Q2, 60-day visit"
$ws.Range("C62").Value = "Anti-S1/S2 IgG Q3, 60-day visit"
$ws.Range("C63").Value = "Anti-S1/S2 IgG Q4, 60-day visit"
$ws.Range("C64").Value = "Ambulatory, acute COVID-19"
$ws.Range("C65").Value = "Hospitalized, acute COVID-19"
$ws.Range("C66").Value = "Oxygen therapy, acute COVID-19"
$ws.Range("C67").Value = "ICU, acute COVID-19"

# --- Remove the hyperlink that was on C6 ("Obesity grade\n@V0") ---
$ws.Hyperlinks.Delete()

# --- Best-effort: restore the view/selection state recorded in the saved file ---
# (final active cell is C64, part of the multi-area selection
#  C40:C52 / C59:C63 / C64 recorded in the original edit)
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("C64").Select() | Out-Null
